$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the Fitness column (C) values for rows 2 through 43 to 7293
$ws.Range("C2:C43").Value = 7293
